# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $ws.Range("F2").Value = 7025
        $ws.Range("F4").Value = 461
        $ws.Range("F23").Value = 2257
        $ws.Range("F25").Value = 255
        $ws.Range("F32").Value = 250
    }
    elseif ($name -eq "全部类型") {
        $ws.Range("F2").Value = 7025
        $ws.Range("F4").Value = 461
        $ws.Range("F24").Value = 2257
        $ws.Range("F26").Value = 255
        $ws.Range("F33").Value = 250
    }
}
